$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: average_doctor / average_doctor_old columns swapped
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Harvard case classification: refreshed "_old" baseline stats (and dependent average_doctor columns)

# Row 4
$ws.Range("E4").Value = 0.473
$ws.Range("F4").Value = 0.056
$ws.Range("G4").Value = 0.236
$ws.Range("N4").Value = 0.482
$ws.Range("O4").Value = 0.062
$ws.Range("P4").Value = 0.25
$ws.Range("Q4").Value = 0.056
$ws.Range("R4").Value = 0.038
$ws.Range("S4").Value = 0.194
$ws.Range("W4").Value = 0.404
$ws.Range("X4").Value = 0.101
$ws.Range("Y4").Value = 0.318
$ws.Range("AI4").Value = 0.404
$ws.Range("AJ4").Value = 0.097
$ws.Range("AK4").Value = 0.311
$ws.Range("AU4").Value = 0.236
$ws.Range("BA4").Value = 2.044
$ws.Range("BB4").Value = 0.152
$ws.Range("BC4").Value = 0.389
$ws.Range("BG4").Value = 0.726
$ws.Range("BH4").Value = 0.149
$ws.Range("BI4").Value = 0.386
$ws.Range("BM4").Value = 0.739
$ws.Range("BN4").Value = 0.065
$ws.Range("BO4").Value = 0.255
$ws.Range("BP4").Value = 0.681
$ws.Range("BQ4").Value = 0.763

# Row 5
$ws.Range("E5").Value = 0.601
$ws.Range("F5").Value = 0.063
$ws.Range("G5").Value = 0.25
$ws.Range("N5").Value = 0.722
$ws.Range("O5").Value = 0.074
$ws.Range("P5").Value = 0.273
$ws.Range("Q5").Value = 0.037
$ws.Range("R5").Value = 0.015
$ws.Range("S5").Value = 0.124
$ws.Range("W5").Value = 0.369
$ws.Range("X5").Value = 0.095
$ws.Range("Y5").Value = 0.308
$ws.Range("AI5").Value = 0.407
$ws.Range("AJ5").Value = 0.094
$ws.Range("AK5").Value = 0.307
$ws.Range("AU5").Value = 0.442
$ws.Range("AV5").Value = 0.07199999999999999
$ws.Range("AW5").Value = 0.268
$ws.Range("BA5").Value = 1.276
$ws.Range("BB5").Value = 0.075
$ws.Range("BC5").Value = 0.274
$ws.Range("BG5").Value = 0.373
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.226
$ws.Range("BM5").Value = 0.512
$ws.Range("BN5").Value = 0.044
$ws.Range("BO5").Value = 0.211
$ws.Range("BP5").Value = 0.425
$ws.Range("BQ5").Value = 0.449

# Row 6
$ws.Range("E6").Value = 0.529
$ws.Range("N6").Value = 0.578
$ws.Range("Q6").Value = 0.045
$ws.Range("W6").Value = 0.386
$ws.Range("AI6").Value = 0.405
$ws.Range("AU6").Value = 0.308
$ws.Range("BA6").Value = 1.565
$ws.Range("BG6").Value = 0.493
$ws.Range("BM6").Value = 0.605
$ws.Range("BP6").Value = 0.522
$ws.Range("BQ6").Value = 0.5620000000000001

# Row 7
$ws.Range("E7").Value = 0.57
$ws.Range("N7").Value = 0.657
$ws.Range("Q7").Value = 0.04
$ws.Range("W7").Value = 0.376
$ws.Range("AI7").Value = 0.406
$ws.Range("AU7").Value = 0.376
$ws.Range("BA7").Value = 1.377
$ws.Range("BG7").Value = 0.413
$ws.Range("BM7").Value = 0.546
$ws.Range("BP7").Value = 0.459
$ws.Range("BQ7").Value = 0.488

# Row 8
$ws.Range("E8").Value = 0.6870000000000001
$ws.Range("F8").Value = 0.073
$ws.Range("G8").Value = 0.27
$ws.Range("N8").Value = 0.8120000000000001
$ws.Range("O8").Value = 0.048
$ws.Range("P8").Value = 0.22
$ws.Range("Q8").Value = 0.041
$ws.Range("S8").Value = 0.163
$ws.Range("W8").Value = 0.446
$ws.Range("X8").Value = 0.113
$ws.Range("Y8").Value = 0.336
$ws.Range("AI8").Value = 0.482
$ws.Range("AJ8").Value = 0.14
$ws.Range("AK8").Value = 0.374
$ws.Range("AU8").Value = 0.39
$ws.Range("AV8").Value = 0.077
$ws.Range("AW8").Value = 0.277
$ws.Range("BA8").Value = 1.761
$ws.Range("BB8").Value = 0.115
$ws.Range("BC8").Value = 0.34
$ws.Range("BG8").Value = 0.5649999999999999
$ws.Range("BH8").Value = 0.11
$ws.Range("BI8").Value = 0.331
$ws.Range("BM8").Value = 0.67
$ws.Range("BN8").Value = 0.06
$ws.Range("BO8").Value = 0.246
$ws.Range("BP8").Value = 0.587
$ws.Range("BQ8").Value = 0.62

# Row 9
$ws.Range("E9").Value = 0.641
$ws.Range("F9").Value = 0.23
$ws.Range("G9").Value = 0.48
$ws.Range("N9").Value = 0.744
$ws.Range("O9").Value = 0.191
$ws.Range("P9").Value = 0.437
$ws.Range("W9").Value = 0.333
$ws.Range("X9").Value = 0.222
$ws.Range("Y9").Value = 0.471
$ws.Range("AI9").Value = 0.436
$ws.Range("AJ9").Value = 0.246
$ws.Range("AK9").Value = 0.496
$ws.Range("BA9").Value = 1.795
$ws.Range("BG9").Value = 0.615
$ws.Range("BH9").Value = 0.237
$ws.Range("BI9").Value = 0.487
$ws.Range("BM9").Value = 0.667
$ws.Range("BN9").Value = 0.222
$ws.Range("BO9").Value = 0.471
$ws.Range("BP9").Value = 0.598
$ws.Range("BQ9").Value = 0.628

# Row 10
$ws.Range("E10").Value = 0.795
$ws.Range("F10").Value = 0.163
$ws.Range("G10").Value = 0.404
$ws.Range("N10").Value = 0.949
$ws.Range("O10").Value = 0.049
$ws.Range("P10").Value = 0.221
$ws.Range("W10").Value = 0.5639999999999999
$ws.Range("X10").Value = 0.246
$ws.Range("Y10").Value = 0.496
$ws.Range("AI10").Value = 0.538
$ws.Range("AU10").Value = 0.41
$ws.Range("AV10").Value = 0.242
$ws.Range("AW10").Value = 0.492
$ws.Range("BA10").Value = 2.18
$ws.Range("BB10").Value = 0.222
$ws.Range("BC10").Value = 0.471
$ws.Range("BG10").Value = 0.6919999999999999
$ws.Range("BH10").Value = 0.213
$ws.Range("BI10").Value = 0.462
$ws.Range("BM10").Value = 0.821
$ws.Range("BN10").Value = 0.147
$ws.Range("BO10").Value = 0.384
$ws.Range("BP10").Value = 0.727
$ws.Range("BQ10").Value = 0.756

# Row 11
$ws.Range("E11").Value = 0.846
$ws.Range("F11").Value = 0.13
$ws.Range("G11").Value = 0.361
$ws.Range("N11").Value = 0.949
$ws.Range("O11").Value = 0.049
$ws.Range("P11").Value = 0.221
$ws.Range("W11").Value = 0.5639999999999999
$ws.Range("X11").Value = 0.246
$ws.Range("Y11").Value = 0.496
$ws.Range("AI11").Value = 0.59
$ws.Range("AJ11").Value = 0.242
$ws.Range("AK11").Value = 0.492
$ws.Range("AU11").Value = 0.5639999999999999
$ws.Range("AV11").Value = 0.246
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.18
$ws.Range("BB11").Value = 0.222
$ws.Range("BC11").Value = 0.471
$ws.Range("BG11").Value = 0.6919999999999999
$ws.Range("BH11").Value = 0.213
$ws.Range("BI11").Value = 0.462
$ws.Range("BM11").Value = 0.821
$ws.Range("BN11").Value = 0.147
$ws.Range("BO11").Value = 0.384
$ws.Range("BP11").Value = 0.727
$ws.Range("BQ11").Value = 0.763

# Row 12
$ws.Range("E12").Value = 1.455
$ws.Range("F12").Value = 0.915
$ws.Range("G12").Value = 0.956
$ws.Range("N12").Value = 1.27
$ws.Range("O12").Value = 0.305
$ws.Range("P12").Value = 0.553
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 0.432
$ws.Range("Y12").Value = 0.657
$ws.Range("AI12").Value = 1.522
$ws.Range("AJ12").Value = 1.293
$ws.Range("AK12").Value = 1.137
$ws.Range("AU12").Value = 2.833
$ws.Range("AV12").Value = 3.306
$ws.Range("AW12").Value = 1.818
$ws.Range("BA12").Value = 3.675
$ws.Range("BG12").Value = 1.148
$ws.Range("BH12").Value = 0.2
$ws.Range("BI12").Value = 0.448
$ws.Range("BM12").Value = 1.219
$ws.Range("BN12").Value = 0.233
$ws.Range("BO12").Value = 0.483
$ws.Range("BP12").Value = 1.225

# Row 13
$ws.Range("E13").Value = 1.432
$ws.Range("F13").Value = 0.307
$ws.Range("G13").Value = 0.554
$ws.Range("N13").Value = 1.73
$ws.Range("O13").Value = 0.449
$ws.Range("P13").Value = 0.67
$ws.Range("W13").Value = 0.985
$ws.Range("X13").Value = 0.194
$ws.Range("Y13").Value = 0.441
$ws.Range("AI13").Value = 1.166
$ws.Range("AJ13").Value = 0.318
$ws.Range("AK13").Value = 0.5639999999999999
$ws.Range("AU13").Value = 2.064
$ws.Range("AV13").Value = 0.333
$ws.Range("AW13").Value = 0.577
$ws.Range("BA13").Value = 2.169
$ws.Range("BB13").Value = 0.303
$ws.Range("BC13").Value = 0.55
$ws.Range("BG13").Value = 0.532
$ws.Range("BH13").Value = 0.049
$ws.Range("BI13").Value = 0.221
$ws.Range("BM13").Value = 0.787
$ws.Range("BN13").Value = 0.173
$ws.Range("BO13").Value = 0.416
$ws.Range("BP13").Value = 0.723
$ws.Range("BQ13").Value = 0.659
